$wb = $excel.ActiveWorkbook

# 1. Fix sheet-name typo: "Fastners" -> "Fasteners".
#    Any formulas referencing the sheet (e.g. Overview!B6 "=Fastners!F29")
#    are automatically rewritten by Excel to use the new name.
$wsFasteners = $wb.Worksheets.Item("Fastners")
$wsFasteners.Name = "Fasteners"

# 2. Move the stored cursor/selection on the "Components" sheet to E39
#    (was A1:H1). Selecting a range on a non-active sheet also makes that
#    sheet the active one, so remember & restore the originally active
#    sheet ("Overview") afterwards to keep the workbook-level active tab
#    unchanged.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsComponents = $wb.Worksheets.Item("Components")
$wsComponents.Range("E39").Select() | Out-Null
$wsOverview.Activate()
